$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.125.77'
$ws.Range('E2').Value = '  +1.13%  '
$ws.Range('D3').Value = '2.376.19'
$ws.Range('E3').Value = '  +3.47%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.12'
$ws.Range('E5').Value = '  +0.75%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.97'
$ws.Range('E6').Value = '  +1.49%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.506'
$ws.Range('E7').Value = '  -0.12%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  +1.29%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.19'
$ws.Range('E10').Value = '  -0.90%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('E12').Value = '  +1.89%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.29'
$ws.Range('E13').Value = '  -4.29%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.80'
$ws.Range('E14').Value = '  +1.53%  '
$ws.Range('D15').Value = '2.748.81'
$ws.Range('E15').Value = '  +3.71%  '
$ws.Range('D16').Value = '2.397.81'
$ws.Range('E16').Value = '  +3.64%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.807'
$ws.Range('E17').Value = '  +3.64%  '
$ws.Range('D18').Value = '43.125.19'
$ws.Range('E18').Value = '  +1.32%  '
$ws.Range('E19').Value = '  -0.89%  '
$ws.Range('E20').Value = '  +4.75%  '
$ws.Range('D21').Value = '0.0₃0886'
$ws.Range('E21').Value = '  -0.25%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.39'
$ws.Range('E22').Value = '  +1.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.40'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.21'
$ws.Range('E24').Value = '  -3.01%  '
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('E25').Value = '  -0.18%  '
$ws.Range('B26').Value = 'PancakeSwap'
$ws.Range('C26').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.43'
$ws.Range('E26').Value = '  +1.03%  '
$ws.Range('E27').Value = '  +2.38%  '
$ws.Range('E28').Value = '  +0.37%  '
$ws.Range('E29').Value = '  +1.15%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '31.36'
$ws.Range('E30').Value = '  -1.13%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  +1.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0748'
$ws.Range('E33').Value = '  +8.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '17.39'
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('E35').Value = '  +5.63%  '
$ws.Range('E36').Value = '  +5.93%  '
$ws.Range('E37').Value = '  -0.54%  '
$ws.Range('E38').Value = '  -2.56%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.79'
$ws.Range('E39').Value = '  +3.95%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.39'
$ws.Range('E40').Value = '  +11.12%  '
$ws.Range('E41').Value = '  -0.12%  '
$ws.Range('D42').Value = '1.960.74'
$ws.Range('E42').Value = '  +0.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '104.07'
$ws.Range('E43').Value = '  -36.91%  '
$ws.Range('E44').Value = '  +0.98%  '
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('E46').Value = '  +0.03%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.13'
$ws.Range('E47').Value = '  -11.13%  '
$ws.Range('B48').Value = 'MultiversX'
$ws.Range('C48').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '52.59'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.50'
$ws.Range('E49').Value = '  +2.02%  '
$ws.Range('B50').Value = 'BitcoinSV'
$ws.Range('C50').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.91'
$ws.Range('E50').Value = '  +1.89%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.14'
$ws.Range('E51').Value = '  +0.96%  '
